$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''29.214.50'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.57%  '
$ws.Range('D3').Value = '''1.827.18'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.82%  '
$ws.Range('D4').Value = '''1.005'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.59%  '
$ws.Range('D5').Value = '''234.15'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.18%  '
$ws.Range('D6').Value = '''0.5951'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.18%  '
$ws.Range('E7').Value = '  +0.47%  '
$ws.Range('D8').Value = '''0.06941'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -6.21%  '
$ws.Range('D9').Value = '''0.2738'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -5.35%  '
$ws.Range('D10').Value = '''23.11'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -7.14%  '
$ws.Range('D11').Value = '''0.07617'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.26%  '
$ws.Range('D12').Value = '''1.843.21'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.17%  '
$ws.Range('D13').Value = '''4.745'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.52%  '
$ws.Range('D14').Value = '''0.6214'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -7.80%  '
$ws.Range('D15').Value = '''0.000009675'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -5.36%  '
$ws.Range('D16').Value = '''78.22'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -4.41%  '
$ws.Range('D17').Value = '''28.877.15'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.78%  '
$ws.Range('D18').Value = '''5.690'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -9.17%  '
$ws.Range('D19').Value = '''220.85'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -5.75%  '
$ws.Range('E20').Value = '  +0.41%  '
$ws.Range('D21').Value = '''11.50'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -6.72%  '
$ws.Range('D22').Value = '''6.840'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -6.38%  '
$ws.Range('D23').Value = '''1.007'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.40%  '
$ws.Range('D24').Value = '''155.69'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.24%  '
$ws.Range('D25').Value = '''7.909'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -6.99%  '
$ws.Range('D26').Value = '''0.1285'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.48%  '
$ws.Range('D27').Value = '''16.44'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.09%  '
$ws.Range('D28').Value = '''0.06643'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -8.53%  '
$ws.Range('D29').Value = '''1.438'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.42%  '
$ws.Range('D30').Value = '''1.438'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.78%  '
$ws.Range('D31').Value = '''3.821'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -5.33%  '
$ws.Range('D32').Value = '''3.739'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -7.43%  '
$ws.Range('D33').Value = '''1.088'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -5.01%  '
$ws.Range('D34').Value = '''1.705'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -6.26%  '
$ws.Range('D35').Value = '''0.6383'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -8.31%  '
$ws.Range('D36').Value = '''2.545'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.11%  '
$ws.Range('D37').Value = '''2.737'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.39%  '
$ws.Range('D38').Value = '''1.182.22'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.17%  '
$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').Value = '''6.501'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.54%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = '''0.01723'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.87%  '
$ws.Range('D41').Value = '''0.8992'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -5.18%  '
$ws.Range('E42').Value = '  +0.43%  '
$ws.Range('D43').Value = '''1.977.84'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.90%  '
$ws.Range('D44').Value = '''100.18'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.84%  '
$ws.Range('D45').Value = '''61.87'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.23%  '
$ws.Range('D46').Value = '''0.00000000114'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.71%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').Value = '''0.05577'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.65%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '''8.457'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.74%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').Value = '''0.4559'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.33%  '
$ws.Range('D50').Value = '''1.563'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -8.21%  '
$ws.Range('D51').Value = '''6.304'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -9.46%  '
